# Assignment1.xlsx update
#
# 1. "PO List" sheet: a handful of literal data cells were updated (the
#    "10storey_next" ranking column K shifted by +1 for most rows, one row's
#    F/H/I/J/K block was reshuffled, and F13 grew from 2 to 3).
# 2. Those data edits change which rows tie for MAX(...) inside the FILTER
#    array formulas on the "Assignment" sheet, so the spilled results there
#    grow from 1/2 cells to 2/3 cells respectively (the formula "ref" widens
#    and the newly spilled cells get the recalculated literal values).
# 3. The active sheet/tab moves from "PO GBW List" to "PO List".

$wb = $excel.ActiveWorkbook

$poList = $wb.Worksheets.Item("PO List")
$assignment = $wb.Worksheets.Item("Assignment")

# --- "PO List" literal data updates -----------------------------------
$poList.Range("K3").Value = 9
$poList.Range("K4").Value = 16
$poList.Range("K5").Value = 19
$poList.Range("K6").Value = 12
$poList.Range("K7").Value = 22
$poList.Range("K8").Value = 18
$poList.Range("K9").Value = 15
$poList.Range("K10").Value = 7
$poList.Range("K11").Value = 14
$poList.Range("K12").Value = 3

$poList.Range("F13").Value = 3
$poList.Range("K13").Value = 2

$poList.Range("K14").Value = 17
$poList.Range("K15").Value = 8
$poList.Range("K16").Value = 11
$poList.Range("K17").Value = 10
$poList.Range("K18").Value = 21
$poList.Range("K19").Value = 6
$poList.Range("K20").Value = 20

$poList.Range("H21").Value = 3
$poList.Range("I21").Value = 3
$poList.Range("J21").Value = 44900
$poList.Range("K21").Value = 1

$poList.Range("K24").Value = 4
$poList.Range("K26").Value = 12
$poList.Range("K29").Value = 5

# --- Recalculate so the FILTER() spill ranges on "Assignment" reflect ---
# --- the data edits above ------------------------------------------------
$wb.Application.Calculate()

# The FILTER array formulas in C10 and C13 now tie with one extra row each
# (Lee Xuan Yen), so their legacy-CSE spill footprint grows by one cell.
# Re-enter each formula over its new, wider range so the stored ref/values
# pick up the additional spilled cell.
$formulaC10 = $assignment.Range("C10").FormulaArray
$assignment.Range("C10:C11").FormulaArray = $formulaC10

$formulaC13 = $assignment.Range("C13").FormulaArray
$assignment.Range("C13:C15").FormulaArray = $formulaC13

$wb.Application.Calculate()

# --- Active sheet moves from "PO GBW List" to "PO List" -----------------
$poList.Activate()
